# Penalty/Reward System update: shift the 16-week forecast window forward
# by one week (each week's start date and MyForecast figure moves to the
# next row) and refresh the dependent Summary-sheet statistics.
#
# NOTE: values are written with a leading apostrophe so Excel stores them
# as literal text (matching the workbook's existing inlineStr / text
# convention) instead of auto-converting date-looking or numeric-looking
# strings into real dates / numbers.

$wb  = $excel.ActiveWorkbook
$wsF = $wb.Worksheets.Item("Forecast Comparison")
$wsS = $wb.Worksheets.Item("Summary")

# ---- Forecast Comparison sheet: Week_Start_Date (B) and MyForecast (D) ----
$weekStarts = @{
    2  = "2025-01-12"
    3  = "2025-01-19"
    4  = "2025-01-26"
    5  = "2025-02-02"
    6  = "2025-02-09"
    7  = "2025-02-16"
    8  = "2025-02-23"
    9  = "2025-03-02"
    10 = "2025-03-09"
    11 = "2025-03-16"
    12 = "2025-03-23"
    13 = "2025-03-30"
    14 = "2025-04-06"
    15 = "2025-04-13"
    16 = "2025-04-20"
    17 = "2025-04-27"
}

$myForecast = @{
    2  = 46
    3  = 49
    4  = 50
    5  = 48
    8  = 32
    9  = 33
    11 = 31
    12 = 31
    13 = 33
    14 = 32
    15 = 30
    16 = 31
}

foreach ($row in $weekStarts.Keys) {
    $wsF.Range("B$row").Value = "'" + $weekStarts[$row]
}

foreach ($row in $myForecast.Keys) {
    $wsF.Range("D$row").Value = $myForecast[$row]
}

# ---- Summary sheet ----
$wsS.Range("B2").Value  = "'2024-02-11 to 2025-01-05"   # Historical Range
$wsS.Range("B4").Value  = "'54"                          # Max Sales
$wsS.Range("B7").Value  = "'15"                           # Std Dev Sales
$wsS.Range("B8").Value  = "'993 units"                    # Total Historical Sales
$wsS.Range("B9").Value  = "'564"                          # Total Forecast (16 Weeks)
$wsS.Range("B10").Value = "'316"                          # Total Forecast (8 Weeks)
$wsS.Range("B11").Value = "'193"                          # Total Forecast (4 Weeks)
$wsS.Range("B12").Value = "'50"                           # Max Forecast
$wsS.Range("B13").Value = "'2025-01-26"                    # Max Forecast Week
$wsS.Range("B15").Value = "'2025-02-09"                    # Min Forecast Week
